# Adds the ability to set a daily (per-day-of-week) lessons amount.
# Expands the sheet from a single "Макс кол-во уроков" column (B) into
# six day-of-week columns (B..G), mirroring the previous max value into
# Monday..Friday (B..F) and adding a Saturday column (G) with either 0
# (grades 5-8) or 5 (grades 9-11) lessons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Класс" -> "Класс/День недели", and day-of-week numbers 1..6
$ws.Range("A1").Value = "Класс/День недели"
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6

# For every data row (2..20), copy the existing "max lessons" value (col B)
# into columns C..F as well (Mon..Fri all use the same previous max), then
# set the Saturday value (col G): grades 5-8 (rows 2-13) have no Saturday
# lessons, grades 9-11 (rows 14-20) have 5.
for ($r = 2; $r -le 20; $r++) {
    $maxLessons = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $maxLessons
    $ws.Cells.Item($r, 4).Value = $maxLessons
    $ws.Cells.Item($r, 5).Value = $maxLessons
    $ws.Cells.Item($r, 6).Value = $maxLessons

    if ($r -le 13) {
        $ws.Cells.Item($r, 7).Value = 0
    } else {
        $ws.Cells.Item($r, 7).Value = 5
    }
}

# Resize column A to fit the new, longer header text.
$ws.Columns.Item(1).ColumnWidth = 15.75

# Adjust the view: scroll so row 2 is the top-left visible row, and
# move the active selection to G21 (just past the data, matching the
# author's final cursor position).
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("G21").Select()
